$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) and "Valor Mora" (column F) rows for
# rows 16..40 are reordered from oldest-to-newest to newest-to-oldest.
# Column B, C, D, G stay identical for every row, only E/F order flips.

$periods = @("1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110")
$values  = @(17667,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,27604)

# Reverse to get the new order (newest period first)
$newPeriods = @($periods[($periods.Length - 1)..0])
$newValues  = @($values[($values.Length - 1)..0])

for ($i = 0; $i -lt $newPeriods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $newPeriods[$i]
    $ws.Cells.Item($row, 6).Value = $newValues[$i]
}
